$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "DinEb858"
$ws.Range("B2").Value = 231009312
$ws.Range("C2").Value = "mspthlv47"
$ws.Range("D2").Value = "TBa73y#$"
$ws.Range("F2").Value = "PNpMBuFk"
$ws.Range("G2").Value = "mMDY"
